$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.959399999999991
$ws.Range("D4").Value = -7.020999999999995
$ws.Range("A9").Value = -20.58679999999997
$ws.Range("D10").Value = -7.67099999999999
$ws.Range("A18").Value = -22.90320000000001
$ws.Range("A20").Value = -22.09900000000003
$ws.Range("C21").Value = -13.41900000000001
